$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Catálogos"

# Rename / re-point the workbook-scoped defined name
$n = $wb.Names.Item(1)
$n.Name = "Catalogos"
$n.RefersTo = "=Catálogos!`$A`$4:`$D`$5"

# Update header row (row 3)
$ws.Range("C3").Value = "Departamento"
$ws.Range("D3").Value = "Activo"
$ws.Range("E3").Clear()

# Update data/placeholder row (row 4)
$ws.Range("C4").Value = "{{item.Departamento}}"
$ws.Range("D4").Value = "{{item.Activo}}"
$ws.Range("E4").Clear()

# Update the selected cell in the sheet view
$ws.Range("G1").Select()
